# Update gh-pages to output generated at 456a3b4
# Refreshes "想去人数" (F) and "最低票价" (G) figures across the 4 sheets
# (展览 / 演出 / 本地生活 / 全部类型) to the latest scraped snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 682
$ws.Range("G6").Value = 88
$ws.Range("F7").Value = 1152
$ws.Range("G7").Value = 70
$ws.Range("F9").Value = 50
$ws.Range("F10").Value = 2041
$ws.Range("F11").Value = 55
$ws.Range("F13").Value = 32
$ws.Range("F16").Value = 1479
$ws.Range("F18").Value = 553
$ws.Range("F19").Value = 392
$ws.Range("F20").Value = 392
$ws.Range("F21").Value = 729
$ws.Range("F22").Value = 452
$ws.Range("F23").Value = 2839
$ws.Range("F24").Value = 399
$ws.Range("F25").Value = 111
$ws.Range("F26").Value = 3208
$ws.Range("F27").Value = 662
$ws.Range("F28").Value = 531
$ws.Range("F29").Value = 234
$ws.Range("F30").Value = 980
$ws.Range("F31").Value = 733
$ws.Range("F32").Value = 41
$ws.Range("F33").Value = 700
$ws.Range("F34").Value = 674

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 147
$ws.Range("F21").Value = 187

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 390

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 682
$ws.Range("G9").Value = 88
$ws.Range("F11").Value = 1152
$ws.Range("G11").Value = 70
$ws.Range("F14").Value = 50
$ws.Range("F15").Value = 390
$ws.Range("F16").Value = 2041
$ws.Range("F17").Value = 2041
$ws.Range("F19").Value = 32
$ws.Range("F27").Value = 1479
$ws.Range("F28").Value = 1479
$ws.Range("F31").Value = 553
$ws.Range("F32").Value = 392
$ws.Range("F33").Value = 392
$ws.Range("F35").Value = 729
$ws.Range("F36").Value = 452
$ws.Range("F38").Value = 2839
$ws.Range("F39").Value = 111
$ws.Range("F40").Value = 3208
$ws.Range("F41").Value = 662
$ws.Range("F42").Value = 531
$ws.Range("F43").Value = 234
$ws.Range("F44").Value = 980
$ws.Range("F49").Value = 733
$ws.Range("F50").Value = 41
$ws.Range("F51").Value = 700
$ws.Range("F52").Value = 674

